$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column I: pixel_size_mm
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true

# Set value for I2
$ws.Range("I2").Value = 1.818

# Update selection to I1
$ws.Range("I1").Select()
